# Resume - Alex Wilber.docx : French text cleanup (Juno -> OLPRODLOC)
#
# Five of the six edits are each the sole run inside their paragraph, so a
# plain Find/Replace is safe. The sixth ("L'art de l'animation...") sits
# between sibling runs that share identical formatting; naively replacing
# its text would make the runtime coalesce it with its neighbours. To keep
# the surrounding runs intact we briefly nudge a reversible character
# property on the single characters immediately touching the target run so
# the formatting no longer matches, perform the text swap, then restore
# those characters' formatting.

$d = $word.ActiveDocument

$d.Content.Find.Execute("CV – Alex Wilber", $true, $false, $false, $false, $false, `
    $true, 1, $false, "CV : Alex Wilber", 2)

$d.Content.Find.Execute("Animation Spark : Concepteur d’animations (Jan 2021 - Présent)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Spark Animation : Concepteur d’animation (depuis janvier 2021)", 2)

$d.Content.Find.Execute("Pixel Studio : Concepteur d’animations (juin 2018 - Déc 2020)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Pixel Studio : Concepteur d’animation (juin 2018 à décembre 2020)", 2)

$d.Content.Find.Execute("Animation flash : Concepteur d’animation junior (sep 2016 - mai 2018)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Flash Animation : Concepteur d’animation junior (septembre 2016 à mai 2018)", 2)

$d.Content.Find.Execute("Master of Arts in Animation, Attend graduation" + [char]160 + ": Dec 2025", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Master en animation, obtention attendue du diplôme" + [char]160 + ": décembre 2025", 2)

# --- the tricky one: preserve the neighbouring " " / "New York : Spark
# --- Press." runs as separate runs instead of letting them merge in.
$old = "L’art de l’animation 3D : Guide pour les débutants."
$new = "The Art of 3D Animation: A Guide for Beginners."

$found = $d.Content
$found.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s = $found.Start
$e = $found.End

$hasPre = $s -gt 0
if ($hasPre) {
    $preChar = $d.Range($s - 1, $s)
    $preChar.Italic = 1
}
$postChar = $d.Range($e, $e + 1)
$postChar.Italic = 1

$target = $d.Range($s, $e)
$target.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

$newEnd = $s + $new.Length
if ($hasPre) {
    $preCharAfter = $d.Range($s - 1, $s)
    $preCharAfter.Italic = 0
}
$postCharAfter = $d.Range($newEnd, $newEnd + 1)
$postCharAfter.Italic = 0

Write-Output "done"
